$d = $word.ActiveDocument

# --- Edit 1: paragraph 1 -----------------------------------------------
# "...and with each new opportunities comes new jobs..."
#   -> "...and with each new opportunity comes new jobs..."
$p1 = $d.Paragraphs(1).Range
$p1.Find.Execute("new opportunities comes new jobs", $false, $false, $false, $false, $false, $true, 1, $false, "new opportunity comes new jobs", 2) | Out-Null

# --- Edit 2: paragraph 2 -------------------------------------------------
# "According to a New Zealand document by aiforum, AI has the opportunity..."
#   -> "According to a New Zealand document by iforum, AI has the opportunity..."
# (only the first "aiforum" mention in this paragraph changes; the later
#  "aiforum"/"Aiforum" mentions in the same paragraph stay as-is, so the
#  search text is scoped tightly to this unique phrase.)
$p2 = $d.Paragraphs(2).Range
$p2.Find.Execute("document by aiforum,", $false, $false, $false, $false, $false, $true, 1, $false, "document by iforum,", 2) | Out-Null
